$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E7 and E8: they were stored as inline-string text ("543237"/"526371"),
# convert them to proper numeric values (543237 / 526371).
$ws.Range("E7").Value = 543237
$ws.Range("E8").Value = 526371

# Append new row 9 (10/06/2024 08:51:31 scrape - MAZDOCK)
$ws.Range("A9").Value = "10/06/2024 08:51:31"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "MAZDOCK"
$ws.Range("D9").Value = "Mazagon Dock Shipbuilders Ltd"

# E9's bsecode must stay TEXT (like the old E7/E8 before they were fixed).
# Build it as text via a helper cell so no stray number format lingers on E9 itself.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "543237"
$helper.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("F9").Value = -1.09
$ws.Range("G9").Value = 3118.75
$ws.Range("H9").Value = 984880

# Append new row 10 (10/06/2024 08:51:31 scrape - NMDC)
$ws.Range("A10").Value = "10/06/2024 08:51:31"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "NMDC"
$ws.Range("D10").Value = "Nmdc Limited"

$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "526371"
$helper.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("F10").Value = -0.97
$ws.Range("G10").Value = 256
$ws.Range("H10").Value = 3945014
